# Row 11 (rule R40) previously showed "R40" in column B; the sheet is
# updated so that cell now holds the text "1" instead.
#
# A plain  $ws.Range("B11").Value = "1"  would be auto-typed as a NUMBER
# (Excel's normal type inference for a numeric-looking literal), but the
# target data is a shared *text* string ("1"), same kind of label as the
# other rule names in that column (R10/R20/R30/...).
#
# To force genuine text without disturbing B11's existing style/number
# format (or allocating a new one), stage a formula that evaluates to the
# text string "1" in a scratch cell, copy it, and paste-special the
# resulting VALUE ONLY into B11 - this keeps B11's own formatting intact
# and just swaps its stored value/type.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$scratch = $ws.Range("Z1")
$scratch.Formula = '="1"'   # formula result is the *text* "1"

$scratch.Copy()
$ws.Range("B11").PasteSpecial(-4163)  # xlPasteValues: value/type only, keep B11's own style

$excel.CutCopyMode = $false
$scratch.Clear()
